$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right=4 (was 5), Wrong=-2 (was -1)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right=80 (was 100), Wrong=-16 (was -8), Max text updated
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -16
$ws.Range("E12").Value = "64 / 112"
